$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 2500
$ws.Range("K64").Value = 2500
$ws.Range("M64").Value = -2252

# Row 67
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 2500
$ws.Range("K67").Value = 2500
$ws.Range("M67").Value = -1642

# Row 74
$ws.Range("H74").Value = 11317.667
$ws.Range("I74").Value = 4476.5
$ws.Range("J74").Value = 25000
$ws.Range("K74").Value = 4476.5
$ws.Range("L74").Value = 25000
$ws.Range("M74").Value = -3540.5
$ws.Range("N74").Value = -26872

# Row 77
$ws.Range("H77").Value = 11317.667
$ws.Range("I77").Value = 4476.5
$ws.Range("J77").Value = 25000
$ws.Range("K77").Value = 22382.5
$ws.Range("L77").Value = 125000
$ws.Range("M77").Value = -17702.5
$ws.Range("N77").Value = -134360

# Row 132
$ws.Range("H132").Value = 4098.591
$ws.Range("I132").Value = 4264.25
$ws.Range("J132").Value = 2442
$ws.Range("K132").Value = 12792.75
$ws.Range("L132").Value = 7326
$ws.Range("M132").Value = -10262.75
$ws.Range("N132").Value = -12386

# Row 139
$ws.Range("H139").Value = 99866.336
$ws.Range("J139").Value = 99866.336
$ws.Range("L139").Value = 99866.336
$ws.Range("N139").Value = -110146.336

# Row 140
$ws.Range("H140").Value = 77578.5
$ws.Range("J140").Value = 77578.5
$ws.Range("L140").Value = 77578.5
$ws.Range("N140").Value = -87938.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1552.0312
$ws.Range("I2").Value = 1741.8636
$ws.Range("J2").Value = 1134.4
$ws.Range("K2").Value = 1741.8636
$ws.Range("L2").Value = 1134.4
$ws.Range("M2").Value = -1628.8636
$ws.Range("N2").Value = -1360.4

# Row 45
$ws.Range("H45").Value = 45456436
$ws.Range("I45").Value = 58824684
$ws.Range("J45").Value = 4392.6
$ws.Range("K45").Value = 58824684
$ws.Range("L45").Value = 4392.6
$ws.Range("M45").Value = -58824307
$ws.Range("N45").Value = -5146.6

# Row 63
$ws.Range("H63").Value = 3289.4736
$ws.Range("I63").Value = 1968.1666
$ws.Range("J63").Value = 5554.5713
$ws.Range("K63").Value = 1968.1666
$ws.Range("L63").Value = 5554.5713
$ws.Range("M63").Value = -1282.1666
$ws.Range("N63").Value = -6926.5713

# Row 66
$ws.Range("H66").Value = 3289.4736
$ws.Range("I66").Value = 1968.1666
$ws.Range("J66").Value = 5554.5713
$ws.Range("K66").Value = 9840.833000000001
$ws.Range("L66").Value = 27772.8565
$ws.Range("M66").Value = -6408.833000000001
$ws.Range("N66").Value = -34636.85649999999

# Row 74
$ws.Range("H74").Value = 21488.732
$ws.Range("I74").Value = 1309.129
$ws.Range("J74").Value = 66172.14
$ws.Range("K74").Value = 1309.129
$ws.Range("L74").Value = 66172.14
$ws.Range("M74").Value = -435.1289999999999
$ws.Range("N74").Value = -67920.14

# Row 77
$ws.Range("H77").Value = 21488.732
$ws.Range("I77").Value = 1309.129
$ws.Range("J77").Value = 66172.14
$ws.Range("K77").Value = 6545.645
$ws.Range("L77").Value = 330860.7
$ws.Range("M77").Value = -2177.645
$ws.Range("N77").Value = -339596.7

# Row 116
$ws.Range("H116").Value = 1552.0312
$ws.Range("I116").Value = 1741.8636
$ws.Range("J116").Value = 1134.4
$ws.Range("K116").Value = 1741.8636
$ws.Range("L116").Value = 1134.4
$ws.Range("M116").Value = 552.1364000000001
$ws.Range("N116").Value = -5722.4

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1552.0312
$ws.Range("I3").Value = 1741.8636
$ws.Range("J3").Value = 1134.4
$ws.Range("K3").Value = 1741.8636
$ws.Range("L3").Value = 1134.4
$ws.Range("M3").Value = -1627.8636
$ws.Range("N3").Value = -1362.4

# Row 35
$ws.Range("H35").Value = 38114.5
$ws.Range("I35").Value = 2500
$ws.Range("K35").Value = 2500
$ws.Range("M35").Value = -2190

# Row 138
$ws.Range("H138").Value = 95133.336
$ws.Range("J138").Value = 95133.336
$ws.Range("L138").Value = 95133.336
$ws.Range("N138").Value = -105413.336

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 975.875
$ws.Range("I22").Value = 855.6667
$ws.Range("J22").Value = 1003.61536
$ws.Range("K22").Value = 855.6667
$ws.Range("L22").Value = 1003.61536
$ws.Range("M22").Value = -505.6667
$ws.Range("N22").Value = -1703.61536

# Row 31
$ws.Range("H31").Value = 2325.0159
$ws.Range("I31").Value = 1540.1578
$ws.Range("K31").Value = 1540.1578
$ws.Range("M31").Value = -1245.1578

# Row 34
$ws.Range("H34").Value = 2325.0159
$ws.Range("I34").Value = 1540.1578
$ws.Range("K34").Value = 1540.1578
$ws.Range("M34").Value = -1338.1578

# Row 132
$ws.Range("H132").Value = 6026.077
$ws.Range("I132").Value = 4834
$ws.Range("J132").Value = 9999.666999999999
$ws.Range("K132").Value = 14502
$ws.Range("L132").Value = 29999.001
$ws.Range("M132").Value = -11972
$ws.Range("N132").Value = -35059.001

# Row 135
$ws.Range("H135").Value = 67966.664
$ws.Range("J135").Value = 67966.664
$ws.Range("L135").Value = 67966.664
$ws.Range("N135").Value = -78106.664

# Row 138
$ws.Range("H138").Value = 88070
$ws.Range("J138").Value = 88188.89
$ws.Range("L138").Value = 88188.89
$ws.Range("N138").Value = -98468.89

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 24477.889
$ws.Range("I80").Value = 39332.668
$ws.Range("J80").Value = 17050.5
$ws.Range("K80").Value = 117998.004
$ws.Range("L80").Value = 51151.5
$ws.Range("M80").Value = -117062.004
$ws.Range("N80").Value = -53023.5

# Row 83
$ws.Range("H83").Value = 24477.889
$ws.Range("I83").Value = 39332.668
$ws.Range("J83").Value = 17050.5
$ws.Range("K83").Value = 353994.012
$ws.Range("L83").Value = 153454.5
$ws.Range("M83").Value = -349314.012
$ws.Range("N83").Value = -162814.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4428.1113
$ws.Range("J70").Value = 4150.625
$ws.Range("L70").Value = 4150.625
$ws.Range("N70").Value = -4690.625

# Row 73
$ws.Range("H73").Value = 4428.1113
$ws.Range("J73").Value = 4150.625
$ws.Range("L73").Value = 4150.625
$ws.Range("N73").Value = -6022.625

# Row 80
$ws.Range("H80").Value = 2274.25
$ws.Range("I80").Value = 2256.2856
$ws.Range("K80").Value = 2256.2856
$ws.Range("M80").Value = -1258.2856

# Row 83
$ws.Range("H83").Value = 2274.25
$ws.Range("I83").Value = 2256.2856
$ws.Range("K83").Value = 11281.428
$ws.Range("M83").Value = -6289.428

# Row 102
$ws.Range("H102").Value = 52020.87
$ws.Range("I102").Value = 115610.445
$ws.Range("K102").Value = 115610.445
$ws.Range("M102").Value = -113988.445

# Row 122
$ws.Range("H122").Value = 152116.1
$ws.Range("I122").Value = 160043.27
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 480129.8099999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -477679.8099999999
$ws.Range("N122").Value = -9400

# Row 132
$ws.Range("H132").Value = 11901.934
$ws.Range("I132").Value = 11805.667
$ws.Range("J132").Value = 12287
$ws.Range("K132").Value = 35417.001
$ws.Range("L132").Value = 36861
$ws.Range("M132").Value = -32887.001
$ws.Range("N132").Value = -41921

# Row 135
$ws.Range("H135").Value = 95467.19
$ws.Range("J135").Value = 95467.19
$ws.Range("L135").Value = 95467.19
$ws.Range("N135").Value = -105607.19

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 36170000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""

# Row 80
$ws.Range("H80").Value = 70100.25
$ws.Range("J80").Value = 70100.25
$ws.Range("L80").Value = 70100.25
$ws.Range("N80").Value = -72096.25

# Row 83
$ws.Range("H83").Value = 70100.25
$ws.Range("J83").Value = 70100.25
$ws.Range("L83").Value = 210300.75
$ws.Range("N83").Value = -220284.75

# Row 139
$ws.Range("H139").Value = 83322
$ws.Range("J139").Value = 83322
$ws.Range("L139").Value = 83322
$ws.Range("N139").Value = -93602

# Row 141
$ws.Range("H141").Value = 119000
$ws.Range("J141").Value = 119000
$ws.Range("L141").Value = 119000
$ws.Range("N141").Value = -129360
